# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 63, shifting the existing
# rows 63:73 down to 64:74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 63; Excel shifts rows 63:73
# down to 64:74, carrying along their content and formatting.
$ws.Rows("63").Insert()

# Populate the new row 63 with the new weekly record. Columns that are
# constant for this market/product block (A, B, C, E, F, G, H, I, J)
# are copied from the (now shifted) row below.
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 44943
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100103
$ws.Range("H63").Value = "Frutos de hueso (carozo)"
$ws.Range("I63").Value = 100103003
$ws.Range("J63").Value = "Damasco"
$ws.Range("K63").Value = "Modesto"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 400
$ws.Range("N63").Value = 20000
$ws.Range("O63").Value = 21000
$ws.Range("P63").Value = 20500
$ws.Range("Q63").Value = "$/caja 16 kilos"
$ws.Range("R63").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S63").Value = 1281
$ws.Range("T63").Value = 16
